$d = $word.ActiveDocument

# --- Change 1: merge split runs around "it's" in the "Using an ability..." sentence ---
$text1 = "Using an ability at full energy reduces it" + [char]8217 + "s cost by 50% or it" + [char]8217 + "s cost for 5 seconds if it" + [char]8217 + "s a channeled ability."
$d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# --- Change 2: merge split runs around "Syncronize" ---
$text2 = "Syncronize" + [char]8211 + " Damaging a block with a basic attack and a ball within 0.5 sec restores energy and equal to 2x the cost of the click"
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# --- Change 3: merge split runs around "Synchronize's" ---
$text3 = "Synchronize" + [char]8217 + "s creates a buff that increases the effectiveness of synchronize by 100% per stack for 2 seconds."
$d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# --- Change 4: add new paragraphs at the end of the document ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)

$p1 = $d.Paragraphs.Add()
$p2 = $d.Paragraphs.Add()
$p3 = $d.Paragraphs.Add()

$p4 = $d.Paragraphs.Add()
$p4.Range.Text = "Retain portion of stages completed up to 100%  (Records a fraction of a stage only while the effect is active and adds them to the total on rebirth.)"

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter([char]11)
$r.Collapse(0)
$r.InsertAfter([char]11)
$r.Collapse(0)
$r.InsertAfter("Rework to the rebirth system:  Instead of giving a currency that is then used to purchase upgrades (effectively no choice), give upgrade points based on how much total has been achieved.")
$r.Collapse(0)
$r.InsertAfter("  This would give actual choice to upgrades instead of just a reason to force players to")
$r.Collapse(0)
$r.InsertAfter(" cache in their rewards every so often.)")
$r.Collapse(0)
$r.InsertAfter("  Maybe give half points while in the life, and the other half when resetting?")
